$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2150.1875
$ws.Range("J70").Value = 2570.3
$ws.Range("L70").Value = 7710.900000000001
$ws.Range("N70").Value = -8250.900000000001
$ws.Range("H73").Value = 2150.1875
$ws.Range("J73").Value = 2570.3
$ws.Range("L73").Value = 7710.900000000001
$ws.Range("N73").Value = -9582.900000000001
$ws.Range("H98").Value = 8313.727999999999
$ws.Range("I98").Value = 6612.875
$ws.Range("J98").Value = 9285.643
$ws.Range("K98").Value = 6612.875
$ws.Range("L98").Value = 9285.643
$ws.Range("M98").Value = -5114.875
$ws.Range("N98").Value = -12281.643
$ws.Range("H107").Value = 1622.6923
$ws.Range("I107").Value = 1554.091
$ws.Range("K107").Value = 1554.091
$ws.Range("M107").Value = 365.9090000000001
$ws.Range("H122").Value = 8313.727999999999
$ws.Range("I122").Value = 6612.875
$ws.Range("J122").Value = 9285.643
$ws.Range("K122").Value = 19838.625
$ws.Range("L122").Value = 27856.929
$ws.Range("M122").Value = -17388.625
$ws.Range("N122").Value = -32756.929

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 4250
$ws.Range("I15").Value = 4000
$ws.Range("J15").Value = 4333.3335
$ws.Range("K15").Value = 4000
$ws.Range("L15").Value = 4333.3335
$ws.Range("M15").Value = -3650
$ws.Range("N15").Value = -5033.3335
$ws.Range("H61").Value = 1835.1428
$ws.Range("I61").Value = 1835.1428
$ws.Range("K61").Value = 1835.1428
$ws.Range("M61").Value = -1623.1428
$ws.Range("H63").Value = 1130365.2
$ws.Range("I63").Value = 2254731.2
$ws.Range("J63").Value = 5999.1665
$ws.Range("K63").Value = 2254731.2
$ws.Range("L63").Value = 5999.1665
$ws.Range("M63").Value = -2254045.2
$ws.Range("N63").Value = -7371.1665
$ws.Range("H66").Value = 1130365.2
$ws.Range("I66").Value = 2254731.2
$ws.Range("J66").Value = 5999.1665
$ws.Range("K66").Value = 11273656
$ws.Range("L66").Value = 29995.8325
$ws.Range("M66").Value = -11270224
$ws.Range("N66").Value = -36859.8325
$ws.Range("H74").Value = 6524.875
$ws.Range("I74").Value = 6739.8
$ws.Range("J74").Value = 6166.6665
$ws.Range("K74").Value = 6739.8
$ws.Range("L74").Value = 6166.6665
$ws.Range("M74").Value = -5865.8
$ws.Range("N74").Value = -7914.6665
$ws.Range("H77").Value = 6524.875
$ws.Range("I77").Value = 6739.8
$ws.Range("J77").Value = 6166.6665
$ws.Range("K77").Value = 33699
$ws.Range("L77").Value = 30833.3325
$ws.Range("M77").Value = -29331
$ws.Range("N77").Value = -39569.3325
$ws.Range("H132").Value = 2637.0625
$ws.Range("I132").Value = 1182.9166
$ws.Range("K132").Value = 3548.7498
$ws.Range("M132").Value = -1018.7498
$ws.Range("H136").Value = 1835.1428
$ws.Range("I136").Value = 1835.1428
$ws.Range("K136").Value = 5505.428400000001
$ws.Range("M136").Value = -2955.428400000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 11515.25
$ws.Range("I20").Value = 2491.1667
$ws.Range("J20").Value = 16929.7
$ws.Range("K20").Value = 2491.1667
$ws.Range("L20").Value = 16929.7
$ws.Range("M20").Value = -2244.1667
$ws.Range("N20").Value = -17423.7
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H99").Value = 3380
$ws.Range("I99").Value = 1777.7778
$ws.Range("J99").Value = 5783.3335
$ws.Range("K99").Value = 1777.7778
$ws.Range("L99").Value = 5783.3335
$ws.Range("M99").Value = -279.7778000000001
$ws.Range("N99").Value = -8779.333500000001
$ws.Range("H134").Value = 3025.1
$ws.Range("I134").Value = 1653.8572
$ws.Range("K134").Value = 4961.571599999999
$ws.Range("M134").Value = -2426.571599999999
$ws.Range("H135").Value = 48791
$ws.Range("J135").Value = 48791
$ws.Range("L135").Value = 48791
$ws.Range("N135").Value = -58931
$ws.Range("H137").Value = 49241
$ws.Range("J137").Value = 49241
$ws.Range("L137").Value = 49241
$ws.Range("N137").Value = -59441

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1895.303
$ws.Range("I58").Value = 1432.5862
$ws.Range("K58").Value = 1432.5862
$ws.Range("M58").Value = -1229.5862
$ws.Range("H99").Value = 10530427
$ws.Range("I99").Value = 28573316
$ws.Range("J99").Value = 5408.3335
$ws.Range("K99").Value = 28573316
$ws.Range("L99").Value = 5408.3335
$ws.Range("M99").Value = -28571818
$ws.Range("N99").Value = -8404.333500000001
$ws.Range("H105").Value = 2225.4285
$ws.Range("I105").Value = 1900
$ws.Range("J105").Value = 2659.3333
$ws.Range("K105").Value = 1900
$ws.Range("L105").Value = 2659.3333
$ws.Range("M105").Value = -153
$ws.Range("N105").Value = -6153.3333
$ws.Range("H126").Value = 10530427
$ws.Range("I126").Value = 28573316
$ws.Range("J126").Value = 5408.3335
$ws.Range("K126").Value = 85719948
$ws.Range("L126").Value = 16225.0005
$ws.Range("M126").Value = -85717478
$ws.Range("N126").Value = -21165.0005
$ws.Range("H132").Value = 5280.143
$ws.Range("I132").Value = 2495.5
$ws.Range("J132").Value = 8993
$ws.Range("K132").Value = 7486.5
$ws.Range("L132").Value = 26979
$ws.Range("M132").Value = -4956.5
$ws.Range("N132").Value = -32039
$ws.Range("H134").Value = 6164
$ws.Range("I134").Value = 6663.1113
$ws.Range("J134").Value = 4666.6665
$ws.Range("K134").Value = 19989.3339
$ws.Range("L134").Value = 13999.9995
$ws.Range("M134").Value = -17454.3339
$ws.Range("N134").Value = -19069.9995
$ws.Range("H136").Value = 1895.303
$ws.Range("I136").Value = 1432.5862
$ws.Range("K136").Value = 4297.7586
$ws.Range("M136").Value = -1747.7586
$ws.Range("H139").Value = 108199.836
$ws.Range("J139").Value = 108199.836
$ws.Range("L139").Value = 108199.836
$ws.Range("N139").Value = -118479.836

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H95").Value = 7333.3335
$ws.Range("J95").Value = 7333.3335
$ws.Range("L95").Value = 22000.0005
$ws.Range("N95").Value = -26118.0005

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 835.7273
$ws.Range("I13").Value = 760.5
$ws.Range("J13").Value = 1036.3334
$ws.Range("K13").Value = 760.5
$ws.Range("L13").Value = 1036.3334
$ws.Range("M13").Value = -621.5
$ws.Range("N13").Value = -1314.3334
$ws.Range("H102").Value = 5208.778
$ws.Range("I102").Value = 3653.8333
$ws.Range("J102").Value = 8318.666999999999
$ws.Range("K102").Value = 3653.8333
$ws.Range("L102").Value = 8318.666999999999
$ws.Range("M102").Value = -2031.8333
$ws.Range("N102").Value = -11562.667
$ws.Range("H132").Value = 4728.091
$ws.Range("I132").Value = 3126.375
$ws.Range("K132").Value = 9379.125
$ws.Range("M132").Value = -6849.125

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 50000
$ws.Range("J24").Value = 50000
$ws.Range("L24").Value = 50000
$ws.Range("N24").Value = -50686
$ws.Range("H122").Value = 11400
$ws.Range("I122").Value = 16500
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 49500
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -47050
$ws.Range("N122").Value = -28900
$ws.Range("H132").Value = 8337.263000000001
$ws.Range("J132").Value = 13810.889
$ws.Range("L132").Value = 41432.667
$ws.Range("N132").Value = -46492.667
$ws.Range("H136").Value = 3746.1428
$ws.Range("I136").Value = 1318.25
$ws.Range("J136").Value = 6983.3335
$ws.Range("K136").Value = 3954.75
$ws.Range("L136").Value = 20950.0005
$ws.Range("M136").Value = -1404.75
$ws.Range("N136").Value = -26050.0005

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 8177
$ws.Range("I24").Value = 3009
$ws.Range("J24").Value = 9899.666999999999
$ws.Range("K24").Value = 3009
$ws.Range("L24").Value = 9899.666999999999
$ws.Range("M24").Value = -2779
$ws.Range("N24").Value = -10359.667
$ws.Range("H132").Value = 17551510
$ws.Range("I132").Value = 10193.909
$ws.Range("J132").Value = 41670816
$ws.Range("K132").Value = 30581.727
$ws.Range("L132").Value = 125012448
$ws.Range("M132").Value = -28051.727
$ws.Range("N132").Value = -125017508
$ws.Range("H136").Value = 4563.0386
$ws.Range("I136").Value = 4170.3125
$ws.Range("J136").Value = 5191.4
$ws.Range("K136").Value = 12510.9375
$ws.Range("L136").Value = 15574.2
$ws.Range("M136").Value = -9960.9375
$ws.Range("N136").Value = -20674.2

Write-Host "All updates applied"
